$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.186238
$ws.Range("H2").Value = 0.558714
$ws.Range("I2").Value = 0.05023668284714279
$ws.Range("J2").Value = 0.05023668284714279
$ws.Range("M2").Value = 1.522526333333333
$ws.Range("N2").Value = 4.567579
$ws.Range("O2").Value = 0.2115373313282365
$ws.Range("P2").Value = 0.2115373313282365
$ws.Range("Q2").Value = 0.2835522592673334
$ws.Range("R2").Value = 2.551970333406
$ws.Range("S2").Value = 0.01062693382426758
$ws.Range("T2").Value = 0.01062693382426758

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.186238
$ws.Range("H3").Value = 0.558714
$ws.Range("I3").Value = 0.05023668284714279
$ws.Range("J3").Value = 0.05023668284714279
$ws.Range("O3").Value = 0.4376697219060474
$ws.Range("P3").Value = 0.4376697219060474
$ws.Range("Q3").Value = 0.5866682617206668
$ws.Range("R3").Value = 5.280014355486001
$ws.Range("S3").Value = 0.02198707501119129
$ws.Range("T3").Value = 0.02198707501119129

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.186238
$ws.Range("H4").Value = 0.558714
$ws.Range("I4").Value = 0.05023668284714279
$ws.Range("J4").Value = 0.05023668284714279
$ws.Range("M4").Value = 2.524809666666667
$ws.Range("N4").Value = 7.574429
$ws.Range("O4").Value = 0.3507929467657161
$ws.Range("P4").Value = 0.3507929467657162
$ws.Range("Q4").Value = 0.4702155027006667
$ws.Range("R4").Value = 4.231939524306
$ws.Range("S4").Value = 0.01762267401168392
$ws.Range("T4").Value = 0.01762267401168393

$ws.Range("I5").Value = 0.659992587420158
$ws.Range("J5").Value = 0.6599925874201579
$ws.Range("M5").Value = 1.522526333333333
$ws.Range("N5").Value = 4.567579
$ws.Range("O5").Value = 0.2115373313282365
$ws.Range("P5").Value = 0.2115373313282365
$ws.Range("Q5").Value = 3.725213900609334
$ws.Range("R5").Value = 33.526925105484
$ws.Range("S5").Value = 0.1396130706392781
$ws.Range("T5").Value = 0.139613070639278

$ws.Range("I6").Value = 0.659992587420158
$ws.Range("J6").Value = 0.6599925874201579
$ws.Range("O6").Value = 0.4376697219060474
$ws.Range("P6").Value = 0.4376697219060474
$ws.Range("S6").Value = 0.2888587721962332
$ws.Range("T6").Value = 0.2888587721962332

$ws.Range("I7").Value = 0.659992587420158
$ws.Range("J7").Value = 0.6599925874201579
$ws.Range("M7").Value = 2.524809666666667
$ws.Range("N7").Value = 7.574429
$ws.Range("O7").Value = 0.3507929467657161
$ws.Range("P7").Value = 0.3507929467657162
$ws.Range("Q7").Value = 6.177532605342668
$ws.Range("R7").Value = 55.597793448084
$ws.Range("S7").Value = 0.2315207445846467
$ws.Range("T7").Value = 0.2315207445846467

$ws.Range("G8").Value = 1.074241333333333
$ws.Range("H8").Value = 3.222724
$ws.Range("I8").Value = 0.2897707297326994
$ws.Range("J8").Value = 0.2897707297326994
$ws.Range("M8").Value = 1.522526333333333
$ws.Range("N8").Value = 4.567579
$ws.Range("O8").Value = 0.2115373313282365
$ws.Range("P8").Value = 0.2115373313282365
$ws.Range("Q8").Value = 1.635560718355111
$ws.Range("R8").Value = 14.720046465196
$ws.Range("S8").Value = 0.06129732686469089
$ws.Range("T8").Value = 0.06129732686469089

$ws.Range("G9").Value = 1.074241333333333
$ws.Range("H9").Value = 3.222724
$ws.Range("I9").Value = 0.2897707297326994
$ws.Range("J9").Value = 0.2897707297326994
$ws.Range("O9").Value = 0.4376697219060474
$ws.Range("P9").Value = 0.4376697219060474
$ws.Range("Q9").Value = 3.38396726605289
$ws.Range("R9").Value = 30.45570539447601
$ws.Range("S9").Value = 0.1268238746986229
$ws.Range("T9").Value = 0.1268238746986229

$ws.Range("G10").Value = 1.074241333333333
$ws.Range("H10").Value = 3.222724
$ws.Range("I10").Value = 0.2897707297326994
$ws.Range("J10").Value = 0.2897707297326994
$ws.Range("M10").Value = 2.524809666666667
$ws.Range("N10").Value = 7.574429
$ws.Range("O10").Value = 0.3507929467657161
$ws.Range("P10").Value = 0.3507929467657162
$ws.Range("Q10").Value = 2.712254902732889
$ws.Range("R10").Value = 24.410294124596
$ws.Range("S10").Value = 0.1016495281693855
$ws.Range("T10").Value = 0.1016495281693855
